# Chapter3-styles.docx — "updated MS + word. added tables. first sketch."
#
# The meaningful, author-intentional edit captured by this commit is a
# redefinition of the custom "Table Caption" paragraph style: it is made to
# stop inheriting the italic run formatting from its base style ("Caption"),
# exactly mirroring how the sibling "Image Caption" style already turns
# italics off. In the underlying OOXML this shows up as a new
#   <w:rPr><w:i w:val="0"/></w:rPr>
# block inside the <w:style w:styleId="TableCaption"> definition.
#
# (The surrounding diff also shows the numbering part's internal `w:tmpl`
# cache IDs and the bibliography custom-XML part's `ds:itemID` GUID getting
# re-rolled — those are Word-internal bookkeeping values that Word itself
# regenerates whenever it resaves a document and are not reachable through
# the Word object model, so there is no COM call that legitimately
# reproduces them; we leave those parts untouched rather than guess at
# hard-coding new random-looking IDs.)

$d = $word.ActiveDocument

$tableCaption = $d.Styles("Table Caption")

# Equivalent of unchecking "Italic" on the style's Format > Font dialog —
# this writes <w:rPr><w:i w:val="0"/></w:rPr> into the style definition.
$tableCaption.Font.Italic = $false
